# Apply the commit "changing content of newFile to see working of branch":
#  - keep "Adding rubbish and commiting ." but drop the _GoBack bookmark from
#    that paragraph (it moves to the new last paragraph)
#  - append a blank paragraph
#  - append a paragraph with the text "Adding this line in master "
#  - append a paragraph that holds the (moved) _GoBack bookmark followed by
#    two textWrapping line breaks

$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark; it will be re-created at the new
# location further down.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Collapse a range to the very end of the document body so the inserted XML
# is appended rather than replacing existing content.
$endRng = $d.Content
$endRng.Collapse(0)

$runPr = '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-IN"/></w:rPr>'

$newParagraphs = @"
<w:p>
  <w:pPr>$runPr</w:pPr>
</w:p>
<w:p>
  <w:pPr>$runPr</w:pPr>
  <w:r>$runPr<w:t xml:space="preserve">Adding this line in master </w:t></w:r>
</w:p>
<w:p>
  <w:pPr>$runPr</w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>$runPr<w:br w:type="textWrapping"/></w:r>
  <w:r>$runPr<w:br w:type="textWrapping"/></w:r>
</w:p>
"@

$package = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  "<w:body>$newParagraphs</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

[void]$endRng.InsertXML($package)
